$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.725.60"
$ws.Range("E2").Value = "'  -2.38%  "
$ws.Range("D3").Value = "'2.234.34"
$ws.Range("E3").Value = "'  -2.44%  "
$ws.Range("E4").Value = "'  +0.33%  "
$ws.Range("D5").Value = "'111.81"
$ws.Range("E5").Value = "'  -9.30%  "
$ws.Range("D6").Value = "'295.51"
$ws.Range("E6").Value = "'  +10.44%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "'  -2.35%  "
$ws.Range("E8").Value = "'  +0.29%  "
$ws.Range("D9").Value = "'0.614"
$ws.Range("E9").Value = "'  -1.05%  "
$ws.Range("D10").Value = "'45.34"
$ws.Range("E10").Value = "'  -7.29%  "
$ws.Range("D11").Value = "'0.0922"
$ws.Range("E11").Value = "'  -2.21%  "
$ws.Range("D12").Value = "'55.19"
$ws.Range("E12").Value = "'  +1.10%  "
$ws.Range("D13").Value = "'8.93"
$ws.Range("E13").Value = "'  -2.84%  "
$ws.Range("E14").Value = "'  -3.52%  "
$ws.Range("D15").Value = "'0.903"
$ws.Range("E15").Value = "'  +0.89%  "
$ws.Range("D16").Value = "'15.19"
$ws.Range("E16").Value = "'  -2.70%  "
$ws.Range("D17").Value = "'2.576.14"
$ws.Range("E17").Value = "'  -2.30%  "
$ws.Range("D18").Value = "'2.264.53"
$ws.Range("E18").Value = "'  -1.21%  "
$ws.Range("D19").Value = "'42.617.64"
$ws.Range("E19").Value = "'  -2.59%  "
$ws.Range("D20").Value = "'7.46"
$ws.Range("E20").Value = "'  +5.81%  "
$ws.Range("D21").Value = "'0.0000107"
$ws.Range("E21").Value = "'  -3.03%  "
$ws.Range("D22").Value = "'73.06"
$ws.Range("E22").Value = "'  +0.62%  "
$ws.Range("D23").Value = "'3.53"
$ws.Range("E23").Value = "'  +22.03%  "
$ws.Range("D24").Value = "'2.31"
$ws.Range("E24").Value = "'  -5.73%  "
$ws.Range("D25").Value = "'230.14"
$ws.Range("E25").Value = "'  -2.53%  "
$ws.Range("D26").Value = "'9.35"
$ws.Range("E26").Value = "'  -3.78%  "
$ws.Range("D27").Value = "'11.94"
$ws.Range("E27").Value = "'  +1.15%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "'  -1.58%  "
$ws.Range("D29").Value = "'39.60"
$ws.Range("E29").Value = "'  -7.40%  "
$ws.Range("E30").Value = "'  -1.35%  "
$ws.Range("E31").Value = "'  -3.96%  "
$ws.Range("D32").Value = "'173.89"
$ws.Range("E32").Value = "'  +0.24%  "
$ws.Range("D33").Value = "'21.17"
$ws.Range("E33").Value = "'  -2.55%  "
$ws.Range("D34").Value = "'0.0890"
$ws.Range("E34").Value = "'  -2.58%  "
$ws.Range("D35").Value = "'5.70"
$ws.Range("E35").Value = "'  -1.65%  "
$ws.Range("D36").Value = "'4.96"
$ws.Range("E36").Value = "'  +4.77%  "
$ws.Range("D37").Value = "'4.28"
$ws.Range("E37").Value = "'  +5.84%  "
$ws.Range("D38").Value = "'0.127"
$ws.Range("E38").Value = "'  -2.26%  "
$ws.Range("D39").Value = "'0.0369"
$ws.Range("E39").Value = "'  -3.19%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("E40").Value = "'  -2.71%  "
$ws.Range("D41").Value = "'2.53"
$ws.Range("E41").Value = "'  -1.70%  "
$ws.Range("D42").Value = "'0.237"
$ws.Range("E42").Value = "'  -1.01%  "
$ws.Range("D43").Value = "'71.34"
$ws.Range("E43").Value = "'  -5.55%  "
$ws.Range("D44").Value = "'13.18"
$ws.Range("E44").Value = "'  -8.37%  "
$ws.Range("E45").Value = "'  +0.26%  "
$ws.Range("E46").Value = "'  -4.06%  "
$ws.Range("D47").Value = "'5.56"
$ws.Range("E47").Value = "'  -6.62%  "
$ws.Range("D48").Value = "'1.31"
$ws.Range("E48").Value = "'  +1.97%  "
$ws.Range("D49").Value = "'105.71"
$ws.Range("E49").Value = "'  +3.68%  "
$ws.Range("D50").Value = "'8.63"
$ws.Range("E50").Value = "'  +0.53%  "
$ws.Range("D51").Value = "'0.0985"
$ws.Range("E51").Value = "'  -2.09%  "
